$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2392.276
$ws.Range("I15").Value = 2392.276
$ws.Range("K15").Value = 7176.828
$ws.Range("M15").Value = -7007.828

$ws.Range("H64").Value = 3029.9
$ws.Range("I64").Value = 3050
$ws.Range("K64").Value = 3050
$ws.Range("M64").Value = -2802

$ws.Range("H67").Value = 3029.9
$ws.Range("I67").Value = 3050
$ws.Range("K67").Value = 3050
$ws.Range("M67").Value = -2192

$ws.Range("M70").ClearContents()
$ws.Range("H70").Value = 1239.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1239.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3718.5
$ws.Range("N70").Value = -4258.5

$ws.Range("M73").ClearContents()
$ws.Range("H73").Value = 1239.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1239.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3718.5
$ws.Range("N73").Value = -5590.5

$ws.Range("H129").Value = 1044.4429
$ws.Range("J129").Value = 1146.2069
$ws.Range("L129").Value = 3438.620699999999
$ws.Range("N129").Value = -13438.6207

$ws.Range("H132").Value = 4441.0435
$ws.Range("I132").Value = 4188.409
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 12565.227
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -10035.227
$ws.Range("N132").Value = -35057

$ws.Range("H138").Value = 2864.6614
$ws.Range("I138").Value = 2059.48
$ws.Range("J138").Value = 3408.7026
$ws.Range("K138").Value = 6178.440000000001
$ws.Range("L138").Value = 10226.1078
$ws.Range("M138").Value = -1038.440000000001
$ws.Range("N138").Value = -20506.1078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1352.0588
$ws.Range("I2").Value = 1468.2
$ws.Range("J2").Value = 1186.1428
$ws.Range("K2").Value = 1468.2
$ws.Range("L2").Value = 1186.1428
$ws.Range("M2").Value = -1355.2
$ws.Range("N2").Value = -1412.1428

$ws.Range("H7").Value = 40222.145
$ws.Range("J7").Value = 40222.145
$ws.Range("L7").Value = 40222.145
$ws.Range("N7").Value = -40450.145

$ws.Range("H32").Value = 793179.8
$ws.Range("I32").Value = 859966.5600000001
$ws.Range("K32").Value = 859966.5600000001
$ws.Range("M32").Value = -859679.5600000001

$ws.Range("H45").Value = 1183.2354
$ws.Range("I45").Value = 1008.3333
$ws.Range("K45").Value = 1008.3333
$ws.Range("M45").Value = -631.3333

$ws.Range("H97").Value = 968.381
$ws.Range("I97").Value = 883.125
$ws.Range("K97").Value = 883.125
$ws.Range("M97").Value = -387.125

$ws.Range("H102").Value = 4490
$ws.Range("I102").Value = 3986.6667
$ws.Range("K102").Value = 3986.6667
$ws.Range("M102").Value = -2364.6667

$ws.Range("H116").Value = 1352.0588
$ws.Range("I116").Value = 1468.2
$ws.Range("J116").Value = 1186.1428
$ws.Range("K116").Value = 1468.2
$ws.Range("L116").Value = 1186.1428
$ws.Range("M116").Value = 825.8
$ws.Range("N116").Value = -5774.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1352.0588
$ws.Range("I3").Value = 1468.2
$ws.Range("J3").Value = 1186.1428
$ws.Range("K3").Value = 1468.2
$ws.Range("L3").Value = 1186.1428
$ws.Range("M3").Value = -1354.2
$ws.Range("N3").Value = -1414.1428

$ws.Range("H20").Value = 25001998
$ws.Range("I20").Value = 1784.68
$ws.Range("J20").Value = 66669020
$ws.Range("K20").Value = 1784.68
$ws.Range("L20").Value = 66669020
$ws.Range("M20").Value = -1537.68
$ws.Range("N20").Value = -66669514

$ws.Range("H80").Value = 1250.3684
$ws.Range("J80").Value = 226
$ws.Range("L80").Value = 226
$ws.Range("N80").Value = -2222

$ws.Range("H83").Value = 1250.3684
$ws.Range("J83").Value = 226
$ws.Range("L83").Value = 1130
$ws.Range("N83").Value = -11114

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 668.5
$ws.Range("I16").Value = 466.4
$ws.Range("J16").Value = 1005.3333
$ws.Range("K16").Value = 466.4
$ws.Range("L16").Value = 1005.3333
$ws.Range("M16").Value = -179.4
$ws.Range("N16").Value = -1579.3333

$ws.Range("H31").Value = 8376.105
$ws.Range("I31").Value = 1702.1818
$ws.Range("K31").Value = 1702.1818
$ws.Range("M31").Value = -1407.1818

$ws.Range("H34").Value = 8376.105
$ws.Range("I34").Value = 1702.1818
$ws.Range("K34").Value = 1702.1818
$ws.Range("M34").Value = -1500.1818

$ws.Range("H105").Value = 637
$ws.Range("I105").Value = 637
$ws.Range("K105").Value = 637
$ws.Range("M105").Value = 1110

$ws.Range("H113").Value = 668.5
$ws.Range("I113").Value = 466.4
$ws.Range("J113").Value = 1005.3333
$ws.Range("K113").Value = 466.4
$ws.Range("L113").Value = 1005.3333
$ws.Range("M113").Value = 1703.6
$ws.Range("N113").Value = -5345.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N3").ClearContents()
$ws.Range("H3").Value = 4015
$ws.Range("I3").Value = 4015
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 12045
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -11933

$ws.Range("H17").Value = 335.4
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 477
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 1431
$ws.Range("M17").Value = -731
$ws.Range("N17").Value = -1769

$ws.Range("H34").Value = 10417164
$ws.Range("J34").Value = 11111631
$ws.Range("L34").Value = 33334893
$ws.Range("N34").Value = -33335061

$ws.Range("H39").Value = 2490
$ws.Range("J39").Value = 2490
$ws.Range("L39").Value = 7470
$ws.Range("N39").Value = -8058

$ws.Range("H55").Value = 1669.2307
$ws.Range("J55").Value = 1791.6666
$ws.Range("L55").Value = 5374.9998
$ws.Range("N55").Value = -5728.9998

$ws.Range("H107").Value = 66667076
$ws.Range("J107").Value = 142857650
$ws.Range("L107").Value = 428572950
$ws.Range("N107").Value = -428576790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6494774.5
$ws.Range("I16").Value = 1166.75
$ws.Range("J16").Value = 14287104
$ws.Range("K16").Value = 1166.75
$ws.Range("L16").Value = 14287104
$ws.Range("M16").Value = -996.75
$ws.Range("N16").Value = -14287444

$ws.Range("H136").Value = 10418582
$ws.Range("I136").Value = 1836.6
$ws.Range("J136").Value = 27779824
$ws.Range("K136").Value = 5509.799999999999
$ws.Range("L136").Value = 83339472
$ws.Range("M136").Value = -2959.799999999999
$ws.Range("N136").Value = -83344572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

$ws.Range("H136").Value = 4030.3333
$ws.Range("I136").Value = 3665.4614
$ws.Range("J136").Value = 4979
$ws.Range("K136").Value = 10996.3842
$ws.Range("L136").Value = 14937
$ws.Range("M136").Value = -8446.3842
$ws.Range("N136").Value = -20037
